$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01
$summary.Range("B4").Value = 0.01
$summary.Range("B5").Value = 0.1
$summary.Range("B6").Value = 2
$summary.Range("B7").Value = 1
$summary.Range("B9").Value = 50

# --- Strategy Status sheet updates (row 4 = MarketMaking) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01
$status.Range("D4").Value = 2
$status.Range("E4").Value = 0.01
$status.Range("F4").Value = 0.01
$status.Range("G4").Value = 50

# --- New trade row (Trade #2) data, shared by "All Trades" and "MarketMaking" sheets ---
# Note: Date-like text ("2026-02-17") needs a leading apostrophe so Excel keeps it
# as text instead of auto-converting it to a date serial number.
$rowValues = @(
    2,
    "'2026-02-17",
    "08:04:43",
    "MarketMaking",
    "UP",
    0.03,
    0.04169,
    "CLOSED",
    38.9683,
    0.01,
    100.01,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.22
)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $ws.Cells.Item(3, $i + 1).Value = $rowValues[$i]
    }
}
